$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (error_id / error_name / description) that was
# breaking the Postgres export; the data rows (E1..E10) shift up to
# become rows 1-10.
$ws.Rows(1).Delete() | Out-Null

# Leave the same "whole row" selection state Excel shows right after a
# row deletion.
$ws.Range("A1:XFD1").Select() | Out-Null
